$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# Sheet "Instructions" (sheet1): rebuild the notes list, inserting two new
# notes (min/max record count) and renaming the "Total ZEVs Supplied" rule
# to "ZEVs Supply Forecast". The sheet is protected, so unprotect first.
# -------------------------------------------------------------------------
$wsInstructions = $wb.Worksheets.Item("Instructions")
$wsInstructions.Unprotect()

$wsInstructions.Cells.Item(7, 1).Value = 'Please note that there must be at least one record in the "Forecast Report" sheet'
$wsInstructions.Cells.Item(9, 1).Value = 'Please note that there must be no more than 2000 records in the "Forecast Report" sheet'
$wsInstructions.Cells.Item(10, 1).ClearContents()
$wsInstructions.Cells.Item(11, 1).Value = 'Please note that no field in the "Forecast Report" sheet may be blank'
$wsInstructions.Cells.Item(12, 1).ClearContents()
$wsInstructions.Cells.Item(13, 1).Value = 'Please note that:'
$wsInstructions.Cells.Item(14, 1).Value = '(1) "Model Year" should be a 4 digit integer'
$wsInstructions.Cells.Item(15, 1).Value = '(2) "Make" should be no more than 250 characters'
$wsInstructions.Cells.Item(16, 1).Value = '(3) "Model" should be no more than 250 characters'
$wsInstructions.Cells.Item(17, 1).Value = '(4) "Type" should be exactly one of: BEV, PHEV, FCEV, EREV'
$wsInstructions.Cells.Item(18, 1).Value = '(5) "Range" should be a real number with no more than 2 decimal places'
$wsInstructions.Cells.Item(19, 1).Value = '(6) "ZEV Class" should be a single, uppercase letter'
$wsInstructions.Cells.Item(20, 1).Value = '(7) "Vehicle Class and Interior Volume" should be no more than 250 characters'
$wsInstructions.Cells.Item(21, 1).Value = '(8) "ZEVs Supply Forecast" should be an integer'

$wsInstructions.Columns.Item(1).ColumnWidth = 73.83333333333334

$wsInstructions.Protect()

# -------------------------------------------------------------------------
# Sheet "Forecast Report" (sheet2): insert the "Range (km)" column between
# "Type" and "ZEV Class", rename "Total ZEVs Supplied" -> "ZEVs Supply
# Forecast", and widen the last column now that bestFit no longer applies.
# -------------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("Forecast Report")

$wsForecast.Cells.Item(1, 5).Value = "Range (km)"
$wsForecast.Cells.Item(1, 6).Value = "ZEV Class"
$wsForecast.Cells.Item(1, 7).Value = "Vehicle Class and Interior Volume"
$wsForecast.Cells.Item(1, 8).Value = "ZEVs Supply Forecast"

$wsForecast.Columns.Item(8).ColumnWidth = 18.5

# -------------------------------------------------------------------------
# Sheet "Dropdowns" (sheet3): reorder the "Station wagon" rows (Small before
# Mid-size) to match the renumbered shared strings.
# -------------------------------------------------------------------------
$wsDropdowns = $wb.Worksheets.Item("Dropdowns")

$wsDropdowns.Cells.Item(8, 3).Value = "Station wagon: Small (less than 130 cu. ft.)"
$wsDropdowns.Cells.Item(9, 3).Value = "Station wagon: Mid-size (130–159 cu. ft.)"
